$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 252
$ws1.Range("F6").Value = 382
$ws1.Range("F7").Value = 239
$ws1.Range("F8").Value = 2276
$ws1.Range("F10").Value = 5656
$ws1.Range("F12").Value = 370

# Sheet "全部类型" (all types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 252
$ws4.Range("F7").Value = 382
$ws4.Range("F8").Value = 239
$ws4.Range("F11").Value = 2276
$ws4.Range("F13").Value = 5656
$ws4.Range("F15").Value = 370
